# Append/update the "取得日時" (retrieved timestamp) column on the
# "ランサーズ" sheet so that every data row (2-12) reflects the latest
# scrape time: 2025-10-24 12:49:25 (JST).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-10-24 12:49:25"

for ($row = 2; $row -le 12; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
